$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- 1. Update simple single-value cells (rows 1-6) ---
$tbl.Cell(1, 1).Range.Text = "0M"
$tbl.Cell(2, 1).Range.Text = "0M"
$tbl.Cell(3, 1).Range.Text = "0M"
$tbl.Cell(4, 1).Range.Text = "3477"
$tbl.Cell(5, 1).Range.Text = "0.00001"
$tbl.Cell(6, 1).Range.Text = "0.00087"

# --- 2. Remove the now-redundant rows 7 and 8 (0.00010 / 0.00003) ---
$tbl.Rows.Item(8).Delete()
$tbl.Rows.Item(7).Delete()

# After the deletions, row numbering shifts up by 2 for everything that
# followed row 8. The row that used to be row 10 (0.00014) is now row 8,
# the row that used to be row 11 (0.00015) is now row 9, and the row that
# used to be row 12 (0.16740) is now row 10.
$tbl.Cell(8, 1).Range.Text = "0.00004"
$tbl.Cell(9, 1).Range.Text = "0.00019"
$tbl.Cell(10, 1).Range.Text = "0.00020"

# --- 3. Insert two new rows after the current row 10 (former 0.16740 row) ---
$newRow1 = $tbl.Rows.Add($tbl.Rows.Item(11))
$tbl.Cell(11, 1).Range.Text = "0.00024"

$newRow2 = $tbl.Rows.Add($tbl.Rows.Item(12))
$tbl.Cell(12, 1).Range.Text = "0.47519"

# --- 4. Collapse the three multi-tab summary rows near the end of the
#         table down to a single value each. These were originally the
#         last three rows in the table and remain so after the edits
#         above (net row count is unchanged: -2 then +2).
$totalRows = $tbl.Rows.Count
$tbl.Cell($totalRows - 2, 1).Range.Text = "99.86"
$tbl.Cell($totalRows - 1, 1).Range.Text = "0.48"
$tbl.Cell($totalRows, 1).Range.Text = "332"
